# Weekly update: insert a new price-report row for "Vega Modelo de Temuco -
# Berenjena" just above the existing row 129. This pushes every row from the
# old 129 through 210 down by one (to 130..211), enlarging the used range
# from A1:R210 to A1:R211 - matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 129; Excel shifts rows 129-210 down to 130-211 and
# extends the sheet's dimension accordingly.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with this week's record.
$ws.Range("A129").Value = 10
$ws.Range("B129").Value = "Vega Modelo de Temuco"
$ws.Range("C129").Value = "La Araucanía"
$ws.Range("D129").Value = 44529
$ws.Range("E129").Value = 9
$ws.Range("F129").Value = 100112001
$ws.Range("G129").Value = "Berenjena"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 85
$ws.Range("K129").Value = 11000
$ws.Range("L129").Value = 13000
$ws.Range("M129").Value = 12059
$ws.Range("N129").Value = "$/caja 60 unidades"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 201
$ws.Range("Q129").Value = 60
$ws.Range("R129").Value = "Hortaliza"
